# 添加R6250H0整机与openEuler 20.03 LTS SP4适配结果
# (adds the VASTAI VG1000/SG100 GPU board-card compatibility row to the
# main compatibility worksheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New compatibility entry -> row 43 (right after the current last row, 42).
# Column layout (row 1 headers): A vendorID, B deviceID, C svID, D ssID,
# E architecture, F os, G driverName, H version, I type, J date, K sha256,
# L driverSize, M chipVendor, N boardModel, O chipModel, P item,
# Q downloadLink, R 备注
$ws.Range("A43").Value = "1ec6"
$ws.Range("B43").Value = "0200"
$ws.Range("C43").Value = "1ec6"
$ws.Range("D43").Value = "0063"
$ws.Range("E43").Value = "aarch64"
$ws.Range("F43").Value = "openEuler 22.03 LTS SP2"
$ws.Range("G43").Value = "vastai_pci"
$ws.Range("H43").Value = "CD9440DD8E1973C48A18416"
$ws.Range("I43").Value = "GPU"
$ws.Range("J43").Value = "2024.1.19"
$ws.Range("K43").Value = "df47d49f3418d5a61002739a246470d35238c8bff4d4ee2fa1081c186c95bf57"
$ws.Range("L43").Value = "29.8M"
$ws.Range("M43").Value = "VASTAI"
$ws.Range("N43").Value = "VG1000"
$ws.Range("O43").Value = "SG100"
$ws.Range("Q43").Value = "https://repo.oepkgs.net/openEuler/rpm/openEuler-22.03-LTS-SP2/contrib/drivers/source/Packages/ddk_pack.tar.gz"

# Move/update the active selection to reflect where editing ended up.
[void]$ws.Range("H53").Select()
